$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.890.68'
$ws.Range("D2").Style = $s
$ws.Range("E2").Value = '  -2.16%  '

$s = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.300.38'
$ws.Range("D3").Style = $s
$ws.Range("E3").Value = '  -5.03%  '

$s = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $s

$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.72'
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = '  -0.85%  '

$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.05'
$ws.Range("D6").Style = $s

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  -2.95%  '

$s = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.296.44'
$ws.Range("D9").Style = $s
$ws.Range("E9").Value = '  -5.16%  '

$ws.Range("E10").Value = '  -3.04%  '

$s = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.58'
$ws.Range("D11").Style = $s
$ws.Range("E11").Value = '  -2.72%  '

$ws.Range("E12").Value = '  +0.93%  '

$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.337'
$ws.Range("D13").Style = $s
$ws.Range("E13").Value = '  -5.07%  '

$s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.08'
$ws.Range("D14").Style = $s
$ws.Range("E14").Value = '  -2.95%  '

$s = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.707.85'
$ws.Range("D15").Style = $s
$ws.Range("E15").Value = '  -5.17%  '

$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.810.64'
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = '  -2.14%  '

$ws.Range("E17").Value = '  -3.26%  '

$s = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.249.93'
$ws.Range("D18").Style = $s
$ws.Range("E18").Value = '  -7.37%  '

$s = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.70'
$ws.Range("D19").Style = $s
$ws.Range("E19").Value = '  -5.03%  '

$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '315.71'
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = '  -3.83%  '

$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.50'
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = '  -4.14%  '

$ws.Range("E23").Value = '  +0.14%  '

$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.16'
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = '  -3.21%  '

$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.172'
$ws.Range("D25").Style = $s
$ws.Range("E25").Value = '  -2.84%  '

$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").Value = '  -6.41%  '

$ws.Range("E28").Value = '  -7.52%  '

$ws.Range("E29").Value = '  -0.73%  '

$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.75'
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = '  -0.38%  '

$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0733'
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = '  -5.45%  '

$ws.Range("E32").Value = '  +3.67%  '

$ws.Range("E33").Value = '  -4.76%  '

$s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.384'
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = '  -4.26%  '

$s = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.81'
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = '  -4.02%  '

$ws.Range("E37").Value = '  -0.12%  '

$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.26'
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = '  -6.37%  '

$s = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.99'
$ws.Range("D39").Style = $s
$ws.Range("E39").Value = '  -5.53%  '

$ws.Range("E40").Value = '  -2.33%  '

$ws.Range("E41").Value = '  -4.95%  '

$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '305.75'
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = '  -6.04%  '

$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.01'
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = '  -3.40%  '

$ws.Range("E44").Value = '  -5.34%  '

$ws.Range("E45").Value = '  -0.94%  '

$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0503'
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = '  -2.89%  '

$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.88'
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = '  -4.37%  '

$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.558'
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = '  -3.40%  '

$ws.Range("E49").Value = '  -2.95%  '

$ws.Range("E50").Value = '  -4.38%  '

$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.03'
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = '  -0.17%  '
